# Update the "Descriptif_numerique" sheet (Rev_Chambre / Total_Facture stats)
$wb = $excel.ActiveWorkbook

$wsDescriptif = $wb.Worksheets.Item("Descriptif_numerique")
$wsDescriptif.Range("C3").Value = 516.177
$wsDescriptif.Range("F3").Value = 651.634
$wsDescriptif.Range("C4").Value = 670.7828
$wsDescriptif.Range("F4").Value = 698.2683
$wsDescriptif.Range("C7").Value = 302.005
$wsDescriptif.Range("F7").Value = 414.8026
$wsDescriptif.Range("C8").Value = 856.75
$wsDescriptif.Range("F8").Value = 1014.5591
$wsDescriptif.Range("C9").Value = 6327.23
$wsDescriptif.Range("F9").Value = 6642.919

# Update the "Regression" sheet (coefficient values)
$wsRegression = $wb.Worksheets.Item("Regression")
$wsRegression.Range("B2").Value = 2.349004527793824
$wsRegression.Range("B3").Value = 0.9994391188903399
$wsRegression.Range("B4").Value = 1.082739916204542
$wsRegression.Range("B5").Value = 1.294326892916941

# Update the "Regression_R2" sheet (R2 value)
$wsRegressionR2 = $wb.Worksheets.Item("Regression_R2")
$wsRegressionR2.Range("A2").Value = 0.99981413114426
